# Finished CPU chapter: log the entry for "Added images of Ivy Bridge and
# Kepler architecture. Write chapter about CPU." on the protocol sheet, and
# rename the built-in cell style from the German "Standard" to "Normal"
# (workbook locale was normalized in this commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New protocol entry: date in column A, activity text in column B.
$ws.Range("A66").Value = 41377
$ws.Range("B66").Value = "Added images of Ivy Bridge and Kepler architecture. Write chapter about CPU."

# Move the active selection down to the next empty row, ready for the
# following entry.
$ws.Range("B71").Select() | Out-Null

# Rename the built-in "Standard" cell style to "Normal".
$styles = $wb.Styles
$styles.Item(1).Delete() | Out-Null
$styles.Add("Normal") | Out-Null
